$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card21")

# --- Header row (row 1): add new column O "Serviced by " ---
# Copy N1's formatting (bold, centered, bordered header style) onto O1
$ws.Cells.Item(1, 14).Copy()
$ws.Cells.Item(1, 15).PasteSpecial(-4122)
$ws.Cells.Item(1, 15).Value = "Serviced by "

# Fix N1 text: drop trailing space ("Correction " -> "Correction")
$ws.Cells.Item(1, 14).Value = "Correction"

# --- Data rows 2-12 ---
# N2:N12 were blank cells -> now contain "nan"
# O2:O12 are new blank cells in the new column; materialize them with the
# same (default/no) formatting as the rest of the data rows by copying N's
# cell format first, then restoring N's value.
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 14).Copy()
    $ws.Cells.Item($r, 15).PasteSpecial(-4122)
    $ws.Cells.Item($r, 14).Value = "nan"
}

$excel.CutCopyMode = 0
